$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated StatQuery text (column C, rows 2-4) - replaces the old malformed
# "all_studies" query with the corrected Programs/Studies/Cases/Samples/
# Case Files/Study Files query.
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['English Setter']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# Update the sheet view: zoom level and active selection moved from B2 to B4.
$excel.ActiveWindow.Zoom = 85
$ws.Range("B4").Select()
